$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell text values scraped on Wed Jul 31 19:47:30 UTC 2024 by GitHub Actions.
# D/E columns hold plain text (prices use "." as a thousands separator in
# some rows, so everything is written/kept as text, never as a number).

$ws.Range("D2").Value = "65.662.49"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "3.268.97"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.29"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.635"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.36%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "3.269.63"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.82"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.397"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").Value = "3.848.00"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("E14").Value = "  -2.65%  "
$ws.Range("D15").Value = "65.826.32"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.314.15"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000161"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "420.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.59%  "
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("E26").Value = "  +5.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.503"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000112"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.80%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  -0.68%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "22.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.61%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.51"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "158.09"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.42"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("D39").Value = "2.840.90"
$ws.Range("E39").Value = "  +3.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.30"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.743"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0633"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "310.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.59"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0266"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("E51").Value = "  -0.37%  "
